# Applies the cryptos.xlsx data refresh described by the commit:
# "Updated cryptos list on Wed Apr  5 09:12:42 UTC 2023 with GitHub Actions"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    # Force a text number format first so Excel does not reinterpret
    # numeric-looking strings (e.g. "21.26") as actual numbers.
    $rng.NumberFormat = "@"
    $rng.Value = $text
    # Restore the default "Normal" style so no stray cell style is
    # left behind (matches the original unstyled data cells).
    $rng.Style = "Normal"
}

Set-TextValue "D2" "28.585.79"
Set-TextValue "E2" "  +1.53%  "
Set-TextValue "D3" "1.913.63"
Set-TextValue "E3" "  +4.64%  "
Set-TextValue "E4" "  -0.02%  "
Set-TextValue "D5" "315.23"
Set-TextValue "E5" "  +1.39%  "
Set-TextValue "E6" "  -0.02%  "
Set-TextValue "D7" "0.5183"
Set-TextValue "E7" "  +4.39%  "
Set-TextValue "E8" "  +0.81%  "
Set-TextValue "D9" "0.09689"
Set-TextValue "E9" "  -1.57%  "
Set-TextValue "E10" "  +3.73%  "
Set-TextValue "E11" "  +2.25%  "
Set-TextValue "D12" "6.535"
Set-TextValue "E12" "  +0.96%  "
Set-TextValue "D13" "21.26"
Set-TextValue "E13" "  +2.97%  "
Set-TextValue "D14" "1.906.13"
Set-TextValue "E14" "  +4.80%  "
Set-TextValue "D15" "7.514"
Set-TextValue "E15" "  +2.75%  "
Set-TextValue "E16" "  -0.01%  "
Set-TextValue "D17" "94.82"
Set-TextValue "E17" "  +2.27%  "
Set-TextValue "E18" "  -0.68%  "
Set-TextValue "D19" "0.06650"
Set-TextValue "E19" "  -0.09%  "
Set-TextValue "D20" "18.22"
Set-TextValue "E20" "  +5.54%  "
Set-TextValue "E21" "  -0.02%  "
Set-TextValue "D22" "6.325"
Set-TextValue "E22" "  +5.26%  "
Set-TextValue "D23" "28.676.23"
Set-TextValue "E23" "  +1.68%  "
Set-TextValue "D24" "11.57"
Set-TextValue "E24" "  +1.93%  "
Set-TextValue "D25" "2.313"
Set-TextValue "E25" "  +2.97%  "
Set-TextValue "D26" "2.679"
Set-TextValue "E26" "  +10.29%  "
Set-TextValue "D27" "2.127.30"
Set-TextValue "E27" "  +4.81%  "
Set-TextValue "D28" "21.30"
Set-TextValue "E28" "  +2.16%  "
Set-TextValue "D29" "158.31"
Set-TextValue "E29" "  -0.31%  "
Set-TextValue "D30" "128.98"
Set-TextValue "E30" "  +1.53%  "
Set-TextValue "D31" "1.116"
Set-TextValue "E31" "  +7.16%  "
Set-TextValue "D32" "0.1084"
Set-TextValue "E32" "  +2.48%  "
Set-TextValue "D33" "5.775"
Set-TextValue "E33" "  +2.92%  "
Set-TextValue "D34" "3.636"
Set-TextValue "E34" "  +0.49%  "
Set-TextValue "D35" "9.902"
Set-TextValue "E35" "  +10.02%  "
Set-TextValue "D36" "0.06797"
Set-TextValue "E36" "  +0.78%  "
Set-TextValue "E37" "  +3.53%  "
Set-TextValue "D38" "1.273"
Set-TextValue "E38" "  +7.61%  "
Set-TextValue "E39" "  +3.65%  "
Set-TextValue "D40" "11.85"
Set-TextValue "E40" "  +3.55%  "
Set-TextValue "B41" "InternetComputer(DFINITY)"
Set-TextValue "C41" "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue "D41" "5.114"
Set-TextValue "E41" "  +2.60%  "
Set-TextValue "B42" "TheSandbox"
Set-TextValue "C42" "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextValue "D42" "0.6482"
Set-TextValue "E42" "  +3.80%  "
Set-TextValue "D43" "1.192"
Set-TextValue "E43" "  +0.93%  "
Set-TextValue "D45" "13.63"
Set-TextValue "E45" "  +3.47%  "
Set-TextValue "D46" "0.6120"
Set-TextValue "E46" "  +2.85%  "
Set-TextValue "D47" "3.773"
Set-TextValue "E47" "  +1.78%  "
Set-TextValue "D48" "1.284"
Set-TextValue "E48" "  +0.58%  "
Set-TextValue "D49" "2.039"
Set-TextValue "E49" "  +4.39%  "
Set-TextValue "D50" "125.07"
Set-TextValue "E50" "  +0.57%  "
Set-TextValue "E51" "  +2.25%  "
